# POCs-Sir.docx edit: "tasks to do added"
#
# Summary of changes (see unified diff):
#  - Several multi-run paragraphs (split apart by spell-check proofErr
#    markers) get collapsed back into a single run with the combined text.
#  - The "Codeship,Shippable,..." bullet gains a new run "- nodejs" and the
#    (hidden) _GoBack bookmark, which used to live in its own empty
#    paragraph right after the "Jenkins-mesos-docker" bullet.
#  - That now-vacated bookmark paragraph becomes a bare empty paragraph.
#  - A new bullet "JENKINS-MESOS from Monday" is added right after the
#    second "Jenkins-mesos-docker" bullet (before the now-bare empty para).
#  - Remaining lone proofErr-wrapped single-run paragraphs (Fp / Mesos /
#    codeship in the second list) simply lose their proofErr wrapper.

$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($index)
    $newPara.Range.Text = $text
    # the original (now shifted one slot later) paragraph is deleted whole,
    # including its paragraph mark, so none of its runs/proofErr survive.
    $d.Paragraphs.Item($index + 1).Range.Delete()
}

# --- 1. "AngularJS-Scala(play) app with |mongodb" -> single run ---
Set-ParaText 1 "AngularJS-Scala(play) app with mongodb"

# --- 2. "Jenkins-|mesos|-|docker" (first bullet) -> single run ---
Set-ParaText 2 "Jenkins-mesos-docker"

# paragraph 3 "Functional programming Eric Meijer 13 chapters" unchanged.

# --- 4. "Codeship,...,study" gains a new run + the _GoBack bookmark ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertAfter("- nodejs")

# Append a single placeholder character so we can anchor a *non-collapsed*
# Range exactly at the new end-of-paragraph position (collapsed ranges
# right on a paragraph boundary don't anchor reliably), add the bookmark
# around it, then delete the placeholder - the bookmark tags stay put,
# immediately after the run and right before the paragraph mark.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertAfter("~")
$p4 = $d.Paragraphs.Item(4)
$endPos = $p4.Range.End
$markerRange = $d.Range($endPos - 2, $endPos - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Range($endPos - 2, $endPos - 1).Delete()

# --- 5. "Mesos| videos" -> single run ---
Set-ParaText 5 "Mesos videos"

# paragraph 6 "Apache Spark<->..." unchanged.

# --- 7. "My gulp project |codeship" -> single run ---
Set-ParaText 7 "My gulp project codeship"

# --- 8. "Jenkins-|mesos|-|docker" (second bullet) -> single run ---
Set-ParaText 8 "Jenkins-mesos-docker"

# --- 9. New bullet "JENKINS-MESOS from Monday" right after bullet 8 ---
# InsertParagraphAfter inherits the ListParagraph/numId=1 formatting of
# paragraph 8, which is exactly what the new bullet needs.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "JENKINS-MESOS from Monday"

# --- 10. Drop the (hidden) _GoBack bookmark from its old paragraph ---
# It now lives in paragraph 4 instead; removing it here (rather than
# deleting/recreating the paragraph) leaves that paragraph as the bare,
# no-pPr empty paragraph the diff expects, right after our new bullet.
$d.Bookmarks.Item("_GoBack").Delete()

# --- Second list: drop the lone proofErr wrappers around single-run text ---
# "Fp", "Mesos", "codeship" keep their text but lose the spell-check markers.
# (indices shifted by +1 vs. the original document because a net one extra
# paragraph was inserted earlier: +1 new bullet "JENKINS-MESOS from Monday")
Set-ParaText 13 "Fp"
Set-ParaText 14 "Mesos"
Set-ParaText 16 "codeship"
